$p = $ppt.ActivePresentation

# --- 1. Update the table's style (slide 16) ---
$s = $p.Slides.Item(16)
foreach ($sh in $s.Shapes) {
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{1D624DC6-FD0A-40BD-86FA-1B06DDFBB8E5}")
    }
}

# --- 2. Swap the Integral / Office theme colour schemes ---
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) -- Office Theme palette,
# expressed as VBA-style 0x00BBGGRR long values for ThemeColorScheme.RGB.
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)
$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
